$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "empty" for the newly populated cells (D6, D7, D8, E8, F8)
$ws.Range("D6").Value = "empty"
$ws.Range("D7").Value = "empty"
$ws.Range("D8").Value = "empty"
$ws.Range("E8").Value = "empty"
$ws.Range("F8").Value = "empty"

# Update the active selection/view as recorded in the workbook
$ws.Range("D1").Select()
$excel.ActiveWindow.ScrollColumn = 2
